# google.docx edit: drop the stale "_GoBack" bookmark from the last
# paragraph (Q12) and append the Q20-Q50 / Q1-Q9 question list that follows
# it, ending with a trailing empty paragraph.

$d = $word.ActiveDocument

# --- 1. Remove the hidden "_GoBack" bookmark -------------------------------
# "_GoBack" is a hidden Word bookmark, so it never shows up while iterating
# $d.Bookmarks, but it can still be reached (and deleted) by name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Build the WordprocessingML for the new paragraphs ------------------
# Each row is: question number, question text, special-formatting marker.
#   $null   -> plain "<number><tab><text>" paragraph
#   "q29"   -> number carries a lastRenderedPageBreak + gramStart/gramEnd
#              proofErr wrapper around the number run
#   "q50"   -> "sharding" is wrapped in spellStart/spellEnd proofErr marks
#   "q1"    -> number wrapped in gramStart/gramEnd proofErr marks
#   "q7"    -> number run carries a lastRenderedPageBreak
# NB: every row is prefixed with the unary "," so each row stays a nested
#     array element instead of being unrolled into $questions.
$questions = @(
    ,@("Q20", "Explain the concept of database locking and different types of locks used for concurrency control.", $null)
    ,@("Q21", "Discuss the concept of materialized views in a database and their advantages and limitations.", $null)
    ,@("Q22", "Explain the concept of database denormalization and scenarios where it may be beneficial.", $null)
    ,@("Q23", "What are the different types of database indexes, and how do they impact query performance?", $null)
    ,@("Q24", "Discuss the differences between OLTP (Online Transaction Processing) and OLAP (Online Analytical Processing) databases.", $null)
    ,@("Q25", "Explain the concept of database backup and recovery strategies, including full and incremental backups.", $null)
    ,@("Q26", "What is the role of a database schema in a DBMS and how does it define data structure?", $null)
    ,@("Q27", "Discuss the differences between a primary key and a unique key in a database.", $null)
    ,@("Q28", "Explain the concept of database views and their benefits in terms of security and query simplification.", $null)
    ,@("Q29", "Can you discuss the concept of database triggers and provide examples?", "q29")
    ,@("Q30", "Explain the purpose of SQL (Structured Query Language) in a DBMS.", $null)
    ,@("Q31", "Discuss the concept of database replication and its importance in achieving data redundancy.", $null)
    ,@("Q32", "Explain the concept of database partitioning and its benefits in performance optimization.", $null)
    ,@("Q33", "What are the challenges associated with scaling a database, and how can they be addressed?", $null)
    ,@("Q34", "Explain the concept of database indexing and different types of indexes.", $null)
    ,@("Q35", "Discuss the differences between clustered and non-clustered indexes in a database.", $null)
    ,@("Q36", "What is the role of a database administrator (DBA) in a DBMS environment?", $null)
    ,@("Q37", "Discuss the concept of database security and measures to protect sensitive data.", $null)
    ,@("Q38", "Explain the concept of database locking and different types of locks for concurrency control.", $null)
    ,@("Q39", "Discuss the differences between optimistic and pessimistic concurrency control in a DBMS.", $null)
    ,@("Q40", "What is database normalization, and why is it important in database design?", $null)
    ,@("Q41", "Discuss the concept of database denormalization and when it may be appropriate.", $null)
    ,@("Q42", "Explain database backup and recovery strategies, including full and incremental backups.", $null)
    ,@("Q43", "How does a distributed database system differ from a centralized database system?", $null)
    ,@("Q44", "Discuss the concept of database normalization and its role in reducing data redundancy.", $null)
    ,@("Q45", "What is the purpose of a database management system (DBMS) in an organization's data infrastructure?", $null)
    ,@("Q46", "Explain the concept of database indexing and its impact on query performance in a DBMS.", $null)
    ,@("Q47", "How does a transaction log ensure data durability and recoverability in a database system?", $null)
    ,@("Q48", "Discuss the advantages and disadvantages of using denormalized databases in certain scenarios.", $null)
    ,@("Q49", "What is the role of a database administrator (DBA) in maintaining data security and integrity?", $null)
    ,@("Q50", $null, "q50")
    ,@("Q1", "Can you explain the role of a database management system (DBMS) and how it manages data?", "q1")
    ,@("Q2", "What are the advantages and disadvantages of using a relational database management system (RDBMS)?", $null)
    ,@("Q3", "Explain the concept of normalization in database design and provide an example.", $null)
    ,@("Q4", "What is the difference between a primary key and a foreign key in a database?", $null)
    ,@("Q5", "How does a DBMS enforce referential integrity in a relational database?", $null)
    ,@("Q6", "Discuss the differences between a clustered index and a non-clustered index.", $null)
    ,@("Q7", "What is a database schema, and how does it define data organization?", "q7")
    ,@("Q8", "Explain the concept of database transactions and their impact on data consistency.", $null)
    ,@("Q9", "What are database triggers and how do they work? Provide an example.", $null)
)

function ConvertTo-XmlText($text) {
    return $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

$sb = New-Object System.Text.StringBuilder
foreach ($q in $questions) {
    $num = $q[0]
    $text = $q[1]
    $special = $q[2]

    [void]$sb.Append("<w:p>")

    if ($special -eq "q29") {
        [void]$sb.Append('<w:proofErr w:type="gramStart"/><w:r><w:lastRenderedPageBreak/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $num))
        [void]$sb.Append('</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $text))
        [void]$sb.Append('</w:t></w:r>')
    }
    elseif ($special -eq "q1") {
        [void]$sb.Append('<w:proofErr w:type="gramStart"/><w:r><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $num))
        [void]$sb.Append('</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:tab/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $text))
        [void]$sb.Append('</w:t></w:r>')
    }
    elseif ($special -eq "q7") {
        [void]$sb.Append('<w:r><w:lastRenderedPageBreak/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $num))
        [void]$sb.Append('</w:t></w:r><w:r><w:tab/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $text))
        [void]$sb.Append('</w:t></w:r>')
    }
    elseif ($special -eq "q50") {
        [void]$sb.Append('<w:r><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $num))
        [void]$sb.Append('</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">Explain the concept of database </w:t></w:r>')
        [void]$sb.Append('<w:proofErr w:type="spellStart"/><w:r><w:t>sharding</w:t></w:r><w:proofErr w:type="spellEnd"/>')
        [void]$sb.Append('<w:r><w:t xml:space="preserve"> and its use in horizontal database scaling.</w:t></w:r>')
    }
    else {
        [void]$sb.Append('<w:r><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $num))
        [void]$sb.Append('</w:t></w:r><w:r><w:tab/><w:t>')
        [void]$sb.Append((ConvertTo-XmlText $text))
        [void]$sb.Append('</w:t></w:r>')
    }

    [void]$sb.Append("</w:p>")
}

# Trailing empty paragraph.
[void]$sb.Append("<w:p/>")

$bodyXml = $sb.ToString()

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $bodyXml +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 3. Insert the new paragraphs right after the last paragraph (Q12) -----
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($packageXml)

Write-Output "Inserted $($questions.Count) question paragraphs; paragraph count is now $($d.Paragraphs.Count)."
